# Add Default Data from JSON file
# The "Blood Type" column (I) values had their casing normalized:
# e.g. "A Positive" -> "A positive" (only the leading word stays
# capitalized, matching the values coming from the default JSON data
# source).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$bloodTypes = @{
    "I2"  = "A positive"
    "I3"  = "Rh negative"
    "I4"  = "Rh positive"
    "I5"  = "O positive"
    "I6"  = "B positive"
    "I7"  = "AB negative"
    "I8"  = "B negative"
    "I9"  = "A negative"
    "I10" = "B positive"
    "I11" = "O negative"
    "I12" = "AB positive"
}

foreach ($addr in $bloodTypes.Keys) {
    $ws.Range($addr).Value = $bloodTypes[$addr]
}

$ws.Activate()
$ws.Range("I12").Select()
